$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.01037833333333333
$ws.Range("H2").Value = 0.031135
$ws.Range("I2").Value = 0.02114284782989566
$ws.Range("J2").Value = 0.02114284782989566
$ws.Range("M2").Value = 5.482938999999999
$ws.Range("N2").Value = 16.448817
$ws.Range("O2").Value = 0.1472261722051079
$ws.Range("P2").Value = 0.147226172205108
$ws.Range("Q2").Value = 0.05690376858833333
$ws.Range("R2").Value = 0.512133917295
$ws.Range("S2").Value = 0.003112780555510611
$ws.Range("T2").Value = 0.003112780555510612

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.01037833333333333
$ws.Range("H3").Value = 0.031135
$ws.Range("I3").Value = 0.02114284782989566
$ws.Range("J3").Value = 0.02114284782989566
$ws.Range("O3").Value = 0.5993885906243068
$ws.Range("P3").Value = 0.5993885906243068
$ws.Range("Q3").Value = 0.2316671631444445
$ws.Range("R3").Value = 2.0850044683
$ws.Range("S3").Value = 0.01267278176254534
$ws.Range("T3").Value = 0.01267278176254534

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.01037833333333333
$ws.Range("H4").Value = 0.031135
$ws.Range("I4").Value = 0.02114284782989566
$ws.Range("J4").Value = 0.02114284782989566
$ws.Range("M4").Value = 9.436472999999999
$ws.Range("N4").Value = 28.309419
$ws.Range("O4").Value = 0.2533852371705853
$ws.Range("P4").Value = 0.2533852371705853
$ws.Range("Q4").Value = 0.09793486228500001
$ws.Range("R4").Value = 0.881413760565
$ws.Range("S4").Value = 0.005357285511839707
$ws.Range("T4").Value = 0.005357285511839707

# Row 5
$ws.Range("G5").Value = 0.4265683333333333
$ws.Range("I5").Value = 0.8690094132698448
$ws.Range("J5").Value = 0.8690094132698448
$ws.Range("M5").Value = 5.482938999999999
$ws.Range("N5").Value = 16.448817
$ws.Range("O5").Value = 0.1472261722051079
$ws.Range("P5").Value = 0.147226172205108
$ws.Range("Q5").Value = 2.338848150998333
$ws.Range("R5").Value = 21.04963335898499
$ws.Range("S5").Value = 0.127940929525926
$ws.Range("T5").Value = 0.127940929525926

# Row 6
$ws.Range("G6").Value = 0.4265683333333333
$ws.Range("I6").Value = 0.8690094132698448
$ws.Range("J6").Value = 0.8690094132698448
$ws.Range("O6").Value = 0.5993885906243068
$ws.Range("P6").Value = 0.5993885906243068
$ws.Range("Q6").Value = 9.521940806544443
$ws.Range("R6").Value = 85.69746725889998
$ws.Range("S6").Value = 0.520874327459068
$ws.Range("T6").Value = 0.520874327459068

# Row 7
$ws.Range("G7").Value = 0.4265683333333333
$ws.Range("I7").Value = 0.8690094132698448
$ws.Range("J7").Value = 0.8690094132698448
$ws.Range("M7").Value = 9.436472999999999
$ws.Range("N7").Value = 28.309419
$ws.Range("O7").Value = 0.2533852371705853
$ws.Range("P7").Value = 0.2533852371705853
$ws.Range("Q7").Value = 4.025300560154999
$ws.Range("R7").Value = 36.22770504139499
$ws.Range("S7").Value = 0.2201941562848508
$ws.Range("T7").Value = 0.2201941562848508

# Row 8
$ws.Range("G8").Value = 0.05392066666666667
$ws.Range("H8").Value = 0.161762
$ws.Range("I8").Value = 0.1098477389002595
$ws.Range("J8").Value = 0.1098477389002595
$ws.Range("M8").Value = 5.482938999999999
$ws.Range("N8").Value = 16.448817
$ws.Range("O8").Value = 0.1472261722051079
$ws.Range("P8").Value = 0.147226172205108
$ws.Range("Q8").Value = 0.2956437261726666
$ws.Range("R8").Value = 2.660793535553999
$ws.Range("S8").Value = 0.01617246212367135
$ws.Range("T8").Value = 0.01617246212367135

# Row 9
$ws.Range("G9").Value = 0.05392066666666667
$ws.Range("H9").Value = 0.161762
$ws.Range("I9").Value = 0.1098477389002595
$ws.Range("J9").Value = 0.1098477389002595
$ws.Range("O9").Value = 0.5993885906243068
$ws.Range("P9").Value = 0.5993885906243068
$ws.Range("Q9").Value = 1.203627545995555
$ws.Range("R9").Value = 10.83264791396
$ws.Range("S9").Value = 0.06584148140269341
$ws.Range("T9").Value = 0.06584148140269341

# Row 10
$ws.Range("G10").Value = 0.05392066666666667
$ws.Range("H10").Value = 0.161762
$ws.Range("I10").Value = 0.1098477389002595
$ws.Range("J10").Value = 0.1098477389002595
$ws.Range("M10").Value = 9.436472999999999
$ws.Range("N10").Value = 28.309419
$ws.Range("O10").Value = 0.2533852371705853
$ws.Range("P10").Value = 0.2533852371705853
$ws.Range("Q10").Value = 0.508820915142
$ws.Range("R10").Value = 4.579388236278
$ws.Range("S10").Value = 0.0278337953738948
$ws.Range("T10").Value = 0.0278337953738948
